# Update the two "undo" FAQ questions so they are distinguishable:
# - Songs section (row 17) now asks about undoing a song change
# - Notation section (row 27) now asks about undoing a notation change
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B17").Value = "Is there a way to undo my last song change?"
$ws.Range("B27").Value = "Is there a way to undo my last notation change?"

# Column B needs to be widened to fit the longer question text
# (ColumnWidth is specified in characters; 44.6 renders as width="45.5" in the
# underlying XML, matching the authored workbook)
$ws.Columns.Item(2).ColumnWidth = 44.6

# Update the active selection to match the authored workbook state
$ws.Range("B31").Select()
